# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" data table: country names at a few ranks
# shift because the underlying ranking changed, several numeric rows are
# refreshed with newer figures, and the "last updated" timestamp advances.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp (row 1) -------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 22:35"

# --- Country names whose rank (row) changed ------------------------------
$ws.Range("A83").Value = "Costa de Marfil"
$ws.Range("A84").Value = "Guatemala"
$ws.Range("A85").Value = "Croacia"

$ws.Range("A102").Value = "Guinea-Bisau"
$ws.Range("A103").Value = "Kenia"

$ws.Range("A111").Value = "Niger"
$ws.Range("A112").Value = "Republica de Chipre"

# --- Refreshed numeric data ----------------------------------------------
function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Row 4 - Estados Unidos
Set-Row 4 1614563 21840 379741 1138769 0 1117 96053

# Row 11 - Peru
Set-Row 11 178945 414 158000 12642 0 33 8303

# Row 55 - Argentina
Set-Row 55 8174 286 3873 4289 0 0 12

# Row 62
Set-Row 62 6677 0 1860 4615 0 2 202

# Row 83
Set-Row 83 2301 70 1100 1172 0 0 29

# Row 84
Set-Row 84 2265 132 159 2061 0 2 45

# Row 85
Set-Row 85 2237 3 1978 162 0 1 97

# Row 102
Set-Row 102 1109 20 42 1061 0 0 6

# Row 103
Set-Row 103 1109 80 366 693 0 0 50

# Row 111
Set-Row 111 924 4 753 111 0 2 60

# Row 112
Set-Row 112 923 1 561 345 0 0 17
